$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(184).Insert()

$ws.Range("A184").Value = 9
$ws.Range("B184").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C184").Value = "Metropolitana"
$ws.Range("D184").Value = 44572
$ws.Range("E184").Value = 13
$ws.Range("F184").Value = 100112052
$ws.Range("G184").Value = "Albahaca"
$ws.Range("H184").Value = "Sin especificar"
$ws.Range("I184").Value = "Primera"
$ws.Range("J184").Value = 160
$ws.Range("K184").Value = 3500
$ws.Range("L184").Value = 4000
$ws.Range("M184").Value = 3750
$ws.Range("N184").Value = "`$/docena de matas"
$ws.Range("O184").Value = "Región Metropolitana"
$ws.Range("P184").Value = 625
$ws.Range("Q184").Value = 6
$ws.Range("R184").Value = "Hortaliza"
